$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.691.08'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '1.817.83'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4444'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3824'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.87'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07684'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.35%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.158'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.382'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.643'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.91%  '
$ws.Range("D16").Value = '1.810.79'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06777'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9989'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.389'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = '28.686.91'
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.401'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '
$ws.Range("D29").Value = '2.014.26'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.302'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.969'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.926'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09359'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2299'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06435'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02364'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6748'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.281'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.247'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.453'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9984'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6194'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.825'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.072'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07131'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.166'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.18%  '
